$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 04:38:41"
$ws1.Range("A3").Value = "Total filas: 18"

# Insert a new row at position 13 (shifts old row13.. down by one)
$ws1.Rows.Item(13).Insert()

# Update rows 8-23 with final data
$ws1.Cells.Item(8, 1).Value = "04:38:41"
$ws1.Cells.Item(8, 2).Value = "04:45"
$ws1.Cells.Item(8, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(8, 4).Value = 7
$ws1.Cells.Item(8, 5).Value = "LP1912"
$ws1.Cells.Item(9, 1).Value = "03:45:24"
$ws1.Cells.Item(9, 2).Value = "04:46"
$ws1.Cells.Item(9, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(9, 4).Value = 61
$ws1.Cells.Item(9, 5).Value = "LP1912"
$ws1.Cells.Item(10, 1).Value = "04:38:41"
$ws1.Cells.Item(10, 2).Value = "04:53"
$ws1.Cells.Item(10, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(10, 4).Value = 15
$ws1.Cells.Item(10, 5).Value = "LP1912"
$ws1.Cells.Item(11, 1).Value = "04:38:41"
$ws1.Cells.Item(11, 2).Value = "05:16"
$ws1.Cells.Item(11, 3).Value = "17_ROMERO"
$ws1.Cells.Item(11, 4).Value = 38
$ws1.Cells.Item(11, 5).Value = "LP1912"
$ws1.Cells.Item(12, 1).Value = "04:38:41"
$ws1.Cells.Item(12, 2).Value = "05:22"
$ws1.Cells.Item(12, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(12, 4).Value = 44
$ws1.Cells.Item(12, 5).Value = "LP1912"
$ws1.Cells.Item(13, 1).Value = "04:38:41"
$ws1.Cells.Item(13, 2).Value = "05:34"
$ws1.Cells.Item(13, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(13, 4).Value = 56
$ws1.Cells.Item(13, 5).Value = "LP1912"
$ws1.Cells.Item(14, 1).Value = "03:45:24"
$ws1.Cells.Item(14, 2).Value = "05:36"
$ws1.Cells.Item(14, 3).Value = "14_ABASTO"
$ws1.Cells.Item(14, 4).Value = 111
$ws1.Cells.Item(14, 5).Value = "LP1912"
$ws1.Cells.Item(15, 1).Value = "04:38:41"
$ws1.Cells.Item(15, 2).Value = "05:46"
$ws1.Cells.Item(15, 3).Value = "15_ABASTO"
$ws1.Cells.Item(15, 4).Value = 68
$ws1.Cells.Item(15, 5).Value = "LP1912"
$ws1.Cells.Item(16, 1).Value = "04:38:41"
$ws1.Cells.Item(16, 2).Value = "05:54"
$ws1.Cells.Item(16, 3).Value = "10_OLMOS"
$ws1.Cells.Item(16, 4).Value = 76
$ws1.Cells.Item(16, 5).Value = "LP1912"
$ws1.Cells.Item(17, 1).Value = "04:38:41"
$ws1.Cells.Item(17, 2).Value = "06:04"
$ws1.Cells.Item(17, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(17, 4).Value = 86
$ws1.Cells.Item(17, 5).Value = "LP1912"
$ws1.Cells.Item(18, 1).Value = "04:38:41"
$ws1.Cells.Item(18, 2).Value = "06:11"
$ws1.Cells.Item(18, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(18, 4).Value = 93
$ws1.Cells.Item(18, 5).Value = "LP1912"
$ws1.Cells.Item(19, 1).Value = "04:38:41"
$ws1.Cells.Item(19, 2).Value = "06:14"
$ws1.Cells.Item(19, 3).Value = "225_HARAS DEL SUR"
$ws1.Cells.Item(19, 4).Value = 96
$ws1.Cells.Item(19, 5).Value = "LP1912"
$ws1.Cells.Item(20, 1).Value = "04:38:41"
$ws1.Cells.Item(20, 2).Value = "06:21"
$ws1.Cells.Item(20, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(20, 4).Value = 103
$ws1.Cells.Item(20, 5).Value = "LP1912"
$ws1.Cells.Item(21, 1).Value = "04:38:41"
$ws1.Cells.Item(21, 2).Value = "06:27"
$ws1.Cells.Item(21, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(21, 4).Value = 109
$ws1.Cells.Item(21, 5).Value = "LP1912"
$ws1.Cells.Item(22, 1).Value = "04:38:41"
$ws1.Cells.Item(22, 2).Value = "06:29"
$ws1.Cells.Item(22, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(22, 4).Value = 111
$ws1.Cells.Item(22, 5).Value = "LP1912"
$ws1.Cells.Item(23, 1).Value = "04:38:41"
$ws1.Cells.Item(23, 2).Value = "06:31"
$ws1.Cells.Item(23, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(23, 4).Value = 113
$ws1.Cells.Item(23, 5).Value = "LP1912"

$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 04:38:41"
$ws2.Range("A3").Value = "Total filas: 4"

# Insert two new rows: one at position 7 (push old row7 to row8) and one at end (row9 handled by direct write)
$ws2.Rows.Item(7).Insert()

$ws2.Cells.Item(6, 1).Value = "04:38:41"
$ws2.Cells.Item(6, 2).Value = "04:45"
$ws2.Cells.Item(6, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(6, 4).Value = 7
$ws2.Cells.Item(6, 5).Value = "LP1912"
$ws2.Cells.Item(7, 1).Value = "03:45:24"
$ws2.Cells.Item(7, 2).Value = "04:46"
$ws2.Cells.Item(7, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(7, 4).Value = 61
$ws2.Cells.Item(7, 5).Value = "LP1912"
$ws2.Cells.Item(8, 1).Value = "04:38:41"
$ws2.Cells.Item(8, 2).Value = "05:34"
$ws2.Cells.Item(8, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(8, 4).Value = 56
$ws2.Cells.Item(8, 5).Value = "LP1912"
$ws2.Cells.Item(9, 1).Value = "04:38:41"
$ws2.Cells.Item(9, 2).Value = "06:11"
$ws2.Cells.Item(9, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(9, 4).Value = 93
$ws2.Cells.Item(9, 5).Value = "LP1912"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 04:38:41"
$ws3.Range("A3").Value = "Total filas: 4"

# Insert a new row at position 6 (push old row6 to row7)
$ws3.Rows.Item(6).Insert()

$ws3.Cells.Item(6, 1).Value = "04:38:41"
$ws3.Cells.Item(6, 2).Value = "05:43"
$ws3.Cells.Item(6, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(6, 4).Value = 65
$ws3.Cells.Item(6, 5).Value = "L6173"
$ws3.Cells.Item(7, 1).Value = "03:45:24"
$ws3.Cells.Item(7, 2).Value = "05:44"
$ws3.Cells.Item(7, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(7, 4).Value = 119
$ws3.Cells.Item(7, 5).Value = "L6173"
$ws3.Cells.Item(8, 1).Value = "04:38:41"
$ws3.Cells.Item(8, 2).Value = "06:08"
$ws3.Cells.Item(8, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(8, 4).Value = 90
$ws3.Cells.Item(8, 5).Value = "L6173"
$ws3.Cells.Item(9, 1).Value = "04:38:41"
$ws3.Cells.Item(9, 2).Value = "06:32"
$ws3.Cells.Item(9, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(9, 4).Value = 114
$ws3.Cells.Item(9, 5).Value = "L6203"
